$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS ---
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value = 890.2719999999999
$ws.Range("C6").Value = 25630.4249615041
$ws.Range("C7").Value = 24625.45591388506
$ws.Range("C8").Value = 22519.310322496553
$ws.Range("C12").Value = 20066.96413394419
$ws.Range("C13").Value = 19670.96413394419
$ws.Range("C14").Value = 12938.964133944195
$ws.Range("C15").Value = 11709.420488444193
$ws.Range("C16").Value = 12123.464488444197
$ws.Range("C18").Value = 1000.0
$ws.Range("C20").Value = 251348.60694873414
$ws.Range("C21").Value = 241493.22723790086
$ws.Range("C22").Value = 220838.99457411078
$ws.Range("C26").Value = 196789.69382414376
$ws.Range("C27").Value = 192906.26042414375
$ws.Range("C28").Value = 126887.8926241438
$ws.Range("C29").Value = 114830.18843300123
$ws.Range("C30").Value = 118890.57302560126
$ws.Range("C32").Value = 9806.649999999998

# --- FUSELAGE ---
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 3000.0
$ws.Range("C6").Value = 3234.0
$ws.Range("D6").Value = 7.8
$ws.Range("C7").Value = 2392.0
$ws.Range("D7").Value = -20.266666666666666
$ws.Range("C8").Value = 3067.0
$ws.Range("D8").Value = 2.2333333333333334
$ws.Range("C9").Value = 2751.0
$ws.Range("D9").Value = -8.3
$ws.Range("D10").Value = -16.966666666666665
$ws.Range("D11").Value = 23.266666666666666
$ws.Range("C12").Value = 2938.833333333333
$ws.Range("D12").Value = -2.038888888888883

# --- WING ---
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 2000.0
$ws.Range("C7").Value = 2424.0
$ws.Range("D7").Value = 21.2
$ws.Range("C8").Value = 1821.0
$ws.Range("D8").Value = -8.95
$ws.Range("C9").Value = 2077.0
$ws.Range("D9").Value = 3.85
$ws.Range("D10").Value = 14.95
$ws.Range("C11").Value = 2489.0
$ws.Range("D11").Value = 24.45
$ws.Range("C12").Value = 2247.0
$ws.Range("D12").Value = 12.35
$ws.Range("C13").Value = 1908.1428571428569
$ws.Range("D13").Value = -4.592857142857143

# --- HORIZONTAL TAIL ---
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 600.0
$ws.Range("D7").Value = -57.0
$ws.Range("D8").Value = -75.16666666666667
$ws.Range("C9").Value = 136.0
$ws.Range("D9").Value = -77.33333333333333
$ws.Range("C10").Value = 181.0
$ws.Range("D10").Value = -69.83333333333333

# --- VERTICAL TAIL ---
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 400.0
$ws.Range("D7").Value = -17.5
$ws.Range("C8").Value = 229.0
$ws.Range("D8").Value = -42.75
$ws.Range("C9").Value = 279.5
$ws.Range("D9").Value = -30.124999999999996

# --- NACELLES ---
$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 900.0
$ws.Range("C3").Value = 628.6666666666665
$ws.Range("D3").Value = -76.71604938271604
$ws.Range("D9").Value = -47.777777777777764
$ws.Range("C10").Value = 514.0
$ws.Range("D10").Value = 14.22222222222225
$ws.Range("D11").Value = -56.88888888888888
$ws.Range("C12").Value = 314.33333333333326
$ws.Range("D16").Value = -47.777777777777764
$ws.Range("C17").Value = 514.0
$ws.Range("D17").Value = 14.22222222222225
$ws.Range("D18").Value = -56.88888888888888
$ws.Range("C19").Value = 314.33333333333326

# --- LANDING GEARS ---
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 500.0
$ws.Range("C5").Value = 840.0
$ws.Range("D5").Value = 68.0
$ws.Range("C6").Value = 1031.0
$ws.Range("D6").Value = 106.2
$ws.Range("C7").Value = 1166.0
$ws.Range("D7").Value = 133.2
$ws.Range("C8").Value = 1005.0
$ws.Range("D8").Value = 101.0
$ws.Range("C9").Value = 1010.5
$ws.Range("D9").Value = 102.09999999999997

# --- SYSTEMS ---
$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 2000.0
$ws.Range("D5").Value = 23.85
$ws.Range("D6").Value = 23.84773931306728
